$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "User Stories"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Observations"

# --- User Stories sheet ---
$ws1.Range("A1").ColumnWidth = 93.85602678571429
$ws1.Range("B1").ColumnWidth = 90.85602678571429

$ws1.Range("A1").Value = 'User Story #1:
As a buyer, I want to add products to my cart so that I can check out and pay for my goods'
$ws1.Range("B1").Value = 'Acceptance Criteria: 
Given a person on a product''s page, When they press on a "Add to cart", then the product will be added with it''s respective price to the user''s cart -Accessible from header option- and are able to pay for their goods afterwards.'
$ws1.Range("A2").Value = 'User Story #2:
As a registered user, I want to see my purchase history so I can track my expenses with this vendor'
$ws1.Range("B2").Value = 'Acceptance Criteria: 
Given a person registered and logged in and already did a purchase, when they click on their "Welcome, USERNAME" button, then the user should be redirected to a profile page showing purchase history'
$ws1.Range("A3").Value = 'User Story #3:
As a visitor, I want to browse the catalgoue of products so I can decide if I will make a purchase before registering'
$ws1.Range("B3").Value = 'Acceptance Criteria: 
Given a person a person wants to browse and shop, when they visit "https://www.demoblaze.com/index.html" , then they will see the exact same catalogue of items as the registered users without restrictions'
$a = $ws1.Range("A1"); $a.WrapText = $true; $a.HorizontalAlignment = -4131; $a.VerticalAlignment = -4160
$b = $ws1.Range("B1"); $b.WrapText = $true; $b.VerticalAlignment = -4160
$a = $ws1.Range("A2"); $a.WrapText = $true; $a.HorizontalAlignment = -4131; $a.VerticalAlignment = -4160
$b = $ws1.Range("B2"); $b.WrapText = $true; $b.VerticalAlignment = -4160
$a = $ws1.Range("A3"); $a.WrapText = $true; $a.HorizontalAlignment = -4131; $a.VerticalAlignment = -4160
$b = $ws1.Range("B3"); $b.WrapText = $true; $b.VerticalAlignment = -4160

$ws1.Range("A1:B1").RowHeight = 130.5
$ws1.Range("A2:B2").RowHeight = 99
$ws1.Range("A3:B3").RowHeight = 102.75

# --- Observations sheet ---
$ws2.Range("A1").ColumnWidth = 168.57087053571428
$ws2.Range("A1").Value = 'There is no inherent value in registering '
$ws2.Range("A2").Value = 'There is no purchase history available '
$ws2.Range("A3").Value = 'When checking out, any invalid input is accepted'
$ws2.Range("A4").Value = 'When adding a product, It has to be added 1 by 1, a +- counter would be better here'
$ws2.Range("A5").Value = 'Homepage Categories acts a filter within same url'
$ws2.Range("A6").Value = 'Pressing next or previous after pressing a category merges all product log again'
$ws2.Range("A7").Value = 'No item sort option, Displayed items are not sorted in any way'
$ws2.Range("A8").Value = 'Cart is properly configured if you''re logged in and refresh/log out-log in again/use a different browser/ use incognito and logging in'
$ws2.Range("A9").Value = 'Openning "https://www.demoblaze.com/cart.html" from a different browser or incognito without logging in displays a gigantic cart, probably from all the unregistered users testing the website '
$ws2.Range("A10").Value = 'Cart remembers your addition if you''re logged out and add items, this should be reset on page refresh or session end. '

# --- view/selection state ---
$ws2.Range("A5").Select()
$excel.ActiveWindow.Zoom = 85
$ws1.Range("D1").Select()
$excel.ActiveWindow.Zoom = 85
